$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "57-28=",
    "25+52=",
    "60-17=",
    "43-5=",
    "30+17=",
    "69-53=",
    "19+53=",
    "13+74=",
    "33+47=",
    "37-8=",
    "60-2=",
    "58-10=",
    "87-40=",
    "53-29=",
    "71-38=",
    "63-52=",
    "57-19=",
    "76-13=",
    "98-14=",
    "85-36=",
    "76-53=",
    "49+48=",
    "73-11=",
    "2+7=",
    "83-51=",
    "15+60=",
    "32+17=",
    "94-54=",
    "3+95=",
    "21-3=",
    "75-42=",
    "56-31=",
    "40+46=",
    "94-28=",
    "79-11=",
    "72-21=",
    "8+9=",
    "27+45=",
    "16+0=",
    "33-18=",
    "2+46=",
    "9-1=",
    "91-80=",
    "81-29=",
    "39+13=",
    "49+44=",
    "18+4=",
    "55+41=",
    "99-26=",
    "83-35=",
    "75-20=",
    "38-25=",
    "78-17=",
    "40+6=",
    "13+49=",
    "51-5=",
    "80-51=",
    "63+9=",
    "9-8=",
    "3+79=",
    "65-4=",
    "25-24=",
    "74+10=",
    "30+56=",
    "21+44=",
    "9-1=",
    "88-70=",
    "3+55=",
    "69-14=",
    "30+46=",
    "98-12=",
    "85-80=",
    "10+72=",
    "19-10=",
    "91-63=",
    "97-95=",
    "49+37=",
    "51+46=",
    "65-61=",
    "0+42=",
    "66-41=",
    "80+9=",
    "56-5=",
    "26+10=",
    "1+68=",
    "67-60=",
    "41+22=",
    "36-35=",
    "28+39=",
    "9+62=",
    "89-11=",
    "7+59=",
    "60-29=",
    "85-29=",
    "96-78=",
    "20+52=",
    "44+36=",
    "90-58=",
    "2+86=",
    "20-6="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count

$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i = $i + 1
    }
}

Write-Output ("Updated " + $i + " cells out of " + $newValues.Count + " (table is " + $rows + "x" + $cols + ")")
